# 9.4.1.1.xlsx -- add the 2020 data column (J) to the table, mirroring the
# formatting of the existing 2019 column (I), and move the selection to J3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column I (rows 3-12, the table body) into column J
# so the new column inherits identical styles/borders (same as Excel's
# "copy cell formatting" behaviour used by the original author).
$ws.Range("I3:I12").Copy() | Out-Null
$ws.Range("J3:J12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# Fill in the new 2020 values.
$ws.Range("J4").Value = 2020
$ws.Range("J5").Value = 253.27664777870578
$ws.Range("J7").Value = 93.236077839070575
$ws.Range("J8").Value = 160
$ws.Range("J10").Value = 69
$ws.Range("J11").Value = 48.5
$ws.Range("J12").Value = 22.8

# Match the saved selection state from the diff.
$ws.Range("J3").Select() | Out-Null
